# ---------------------------------------------------------------------------
# Regenerate the SharePoint "content type schema" custom-XML part metadata
# (customXml/item4.xml + its itemProps4.xml companion) with fresh
# versionID / fieldsID / itemID stamps, and drop the now-unused schemaRefs
# list from itemProps4.xml (Word collapses it to a bare, self-closing
# <ds:datastoreItem .../> once nothing references those schema URIs anymore).
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

$oldVersionID = "1cc4cf9d95b5e2d14d7aabb44ca49f5e"
$newVersionID = "ce94caacb4a5cc228342027e3189af2c"
$oldFieldsID  = "a33e6829bf21261855124b7b230b6e9c"
$newFieldsID  = "5f85a36ab557a4a47cd270a1ee4435c0"

$oldItemID = "{8627FC41-1F38-4246-9AAB-AAEA609F6DC1}"
$newItemID = "{3DF06557-A34E-42E3-A45A-4B2F0DF75AF8}"

$contentTypeNs = "http://schemas.microsoft.com/office/2006/metadata/contentType"

# Locate the two related parts: the content-type schema (item4.xml) and its
# itemProps "datastore" companion, identified by namespace / itemID rather
# than by a fixed collection index (SharePoint can renumber the collection).
$schemaPart = $null
$propsPart  = $null
for ($i = 1; $i -le $d.CustomXMLParts.Count; $i++) {
    $part = $d.CustomXMLParts.Item($i)
    try {
        if ($part.NamespaceURI -eq $contentTypeNs) { $schemaPart = $part }
    } catch { }
    try {
        if ($part.XML -like "*$oldItemID*") { $propsPart = $part }
    } catch { }
}
if ($schemaPart -eq $null) {
    try { $schemaPart = $d.CustomXMLParts.SelectByNamespace($contentTypeNs).Item(1) } catch { }
}
if ($schemaPart -eq $null) {
    try { $schemaPart = $d.CustomXMLParts.SelectByID($oldItemID).Item(1) } catch { }
}

# --- customXml/item4.xml: bump ma:versionID and ma:fieldsID -----------------
if ($schemaPart -ne $null) {
    $xml = $schemaPart.XML
    if ([string]::IsNullOrEmpty($xml)) { $xml = $schemaPart.Text }
    if (-not [string]::IsNullOrEmpty($xml)) {
        $updated = $xml.Replace($oldVersionID, $newVersionID).Replace($oldFieldsID, $newFieldsID)

        $applied = $false
        try { $schemaPart.XML = $updated; $applied = $true } catch { }
        if (-not $applied) {
            try { $schemaPart.Text = $updated; $applied = $true } catch { }
        }
        if (-not $applied) {
            try {
                $schemaPart.DocumentElement.SetAttribute("ma:versionID", $newVersionID)
                $root = $schemaPart.SelectSingleNode("//*[local-name()='schema'][@ma:fieldsID]")
                if ($root -ne $null) { $root.SetAttribute("ma:fieldsID", $newFieldsID) }
                $applied = $true
            } catch { }
        }
        if (-not $applied) {
            # Last resort: replace the whole part (delete + re-add with the
            # corrected XML), which is the documented way to "edit" a
            # CustomXMLPart since CustomXMLPart.XML is otherwise read-only.
            try {
                $schemaPart.Delete()
                $d.CustomXMLParts.Add($updated) | Out-Null
            } catch { }
        }
    }
}

# --- customXml/itemProps4.xml: new itemID, drop the schemaRefs list --------
$newPropsXml = '<?xml version="1.0" encoding="utf-8"?>' + "`r`n" + '<ds:datastoreItem xmlns:ds="http://schemas.openxmlformats.org/officeDocument/2006/customXml" ds:itemID="' + $newItemID + '"/>'

if ($propsPart -eq $null -and $schemaPart -ne $null) {
    try { $propsPart = $schemaPart.XMLMapping.CustomXMLPart } catch { }
}

if ($propsPart -ne $null) {
    $applied2 = $false
    try { $propsPart.XML = $newPropsXml; $applied2 = $true } catch { }
    if (-not $applied2) {
        try {
            $propsPart.DocumentElement.SetAttribute("ds:itemID", $newItemID)
            $refs = $propsPart.SelectSingleNode("//*[local-name()='schemaRefs']")
            if ($refs -ne $null) { $refs.ParentNode.RemoveChild($refs) | Out-Null }
            $applied2 = $true
        } catch { }
    }
    if (-not $applied2) {
        try {
            $propsPart.Delete()
            $d.CustomXMLParts.Add($newPropsXml) | Out-Null
        } catch { }
    }
}

$d.Saved = $false
